# Last update from John!
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "minAcro" renamed to "Misc Values" (columns B and D both mirror this label)
$ws.Range("D5").Value = "Misc Values"
$ws.Range("B5").Value = "Misc Values"

# Row 18 was blank; it now documents the new "mag calibration" write/send pair
$ws.Range("D18").Value = "Send mag cal values"
$ws.Range("B18").Value = "Write mag calibration values"

# Newly documented ISR codes that were previously blank
$ws.Range("D20").Value = "Send Pressure Altitude"
$ws.Range("D24").Value = "Send Command in Detent Discretes"

# Row 16: relabel the calibration-values message (now specific to accel calibration)
$ws.Range("D16").Value = "Send accel calibration values"

$ws.Range("D26").Value = "Send 100 Hz loop time"

# New highlighted font (bold, light blue) applied to the accel/mag calibration rows
$highlight = @("B16", "D16", "D17", "B18", "D18")
foreach ($addr in $highlight) {
    $cell = $ws.Range($addr)
    $cell.Font.Bold = $true
    $cell.Font.Color = 15773696
}

# Restore the user's last selection before saving
$ws.Range("I30").Select()
